# Update ODI appearance counts (column C) for specific players, incrementing by 1
# as additional scraped match data added one more ODI appearance for each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 7, 17, 20, 22, 24, 26, 30, 31, 34, 37)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    $cell.Value = $current + 1
}
